{"js": "// Update each three-digit-by-one-digit multiplication fact in the table\n// to the new problem/answer pair from the commit's regenerated output.\nconst replacements = [\n  [\"922\u00d74=3688\", \"922\u00d77=6454\"],\n  [\"702\u00d73=2106\", \"664\u00d74=2656\"],\n  [\"976\u00d76=5856\", \"727\u00d78=5816\"],\n  [\"317\u00d72=634\", \"665\u00d76=3990\"],\n  [\"235\u00d76=1410\", \"405\u00d79=3645\"],\n  [\"356\u00d77=2492\", \"323\u00d76=1938\"],\n  [\"144\u00d77=1008\", \"132\u00d75=660\"],\n  [\"494\u00d72=988\", \"498\u00d74=1992\"],\n  [\"253\u00d79=2277\", \"983\u00d72=1966\"],\n  [\"913\u00d72=1826\", \"318\u00d73=954\"],\n  [\"838\u00d75=4190\", \"999\u00d74=3996\"],\n  [\"188\u00d74=752\", \"231\u00d74=924\"],\n  [\"721\u00d76=4326\", \"663\u00d73=1989\"],\n  [\"426\u00d75=2130\", \"357\u00d72=714\"],\n  [\"587\u00d78=4696\", \"638\u00d74=2552\"],\n  [\"797\u00d79=7173\", \"135\u00d77=945\"],\n  [\"883\u00d78=7064\", \"319\u00d79=2871\"],\n  [\"398\u00d73=1194\", \"326\u00d78=2608\"],\n  [\"709\u00d72=1418\", \"254\u00d77=1778\"],\n  [\"629\u00d76=3774\", \"647\u00d76=3882\"],\n  [\"788\u00d79=7092\", \"572\u00d76=3432\"],\n  [\"450\u00d78=3600\", \"359\u00d79=3231\"],\n  [\"154\u00d75=770\", \"556\u00d77=3892\"],\n  [\"665\u00d72=1330\", \"562\u00d72=1124\"],\n  [\"204\u00d77=1428\", \"333\u00d79=2997\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication fact in the table\n# to the new problem/answer pair from the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"922\u00d74=3688\", \"922\u00d77=6454\"),\n    @(\"702\u00d73=2106\", \"664\u00d74=2656\"),\n    @(\"976\u00d76=5856\", \"727\u00d78=5816\"),\n    @(\"317\u00d72=634\",  \"665\u00d76=3990\"),\n    @(\"235\u00d76=1410\", \"405\u00d79=3645\"),\n    @(\"356\u00d77=2492\", \"323\u00d76=1938\"),\n    @(\"144\u00d77=1008\", \"132\u00d75=660\"),\n    @(\"494\u00d72=988\",  \"498\u00d74=1992\"),\n    @(\"253\u00d79=2277\", \"983\u00d72=1966\"),\n    @(\"913\u00d72=1826\", \"318\u00d73=954\"),\n    @(\"838\u00d75=4190\", \"999\u00d74=3996\"),\n    @(\"188\u00d74=752\",  \"231\u00d74=924\"),\n    @(\"721\u00d76=4326\", \"663\u00d73=1989\"),\n    @(\"426\u00d75=2130\", \"357\u00d72=714\"),\n    @(\"587\u00d78=4696\", \"638\u00d74=2552\"),\n    @(\"797\u00d79=7173\", \"135\u00d77=945\"),\n    @(\"883\u00d78=7064\", \"319\u00d79=2871\"),\n    @(\"398\u00d73=1194\", \"326\u00d78=2608\"),\n    @(\"709\u00d72=1418\", \"254\u00d77=1778\"),\n    @(\"629\u00d76=3774\", \"647\u00d76=3882\"),\n    @(\"788\u00d79=7092\", \"572\u00d76=3432\"),\n    @(\"450\u00d78=3600\", \"359\u00d79=3231\"),\n    @(\"154\u00d75=770\",  \"556\u00d77=3892\"),\n    @(\"665\u00d72=1330\", \"562\u00d72=1124\"),\n    @(\"204\u00d77=1428\", \"333\u00d79=2997\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
